$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44188

# Row 3
$ws.Range("D3").Value = 44210
$ws.Range("J3").Value = 8800
$ws.Range("K3").Value = 2500
$ws.Range("M3").Value = 2750
$ws.Range("O3").Value = "Provincia de Chacabuco"
$ws.Range("P3").Value = 28

# Row 4
$ws.Range("D4").Value = 44214
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 7000
$ws.Range("K4").Value = 3000
$ws.Range("L4").Value = 3000
$ws.Range("M4").Value = 3000
$ws.Range("O4").Value = "Provincia de Chacabuco"
$ws.Range("P4").Value = 30

# Row 5
$ws.Range("D5").Value = 44204
$ws.Range("J5").Value = 7000

# Row 6
$ws.Range("D6").Value = 44245
$ws.Range("J6").Value = 9000
$ws.Range("O6").Value = "Región Metropolitana"

# Row 7
$ws.Range("D7").Value = 44245
$ws.Range("I7").Value = "Segunda"
$ws.Range("J7").Value = 5000
$ws.Range("K7").Value = 2500
$ws.Range("L7").Value = 2500
$ws.Range("M7").Value = 2500
$ws.Range("O7").Value = "Región Metropolitana"
$ws.Range("P7").Value = 25

# Row 8
$ws.Range("D8").Value = 44181
$ws.Range("J8").Value = 12000

# Row 9
$ws.Range("D9").Value = 44229
$ws.Range("J9").Value = 16000

# Row 10
$ws.Range("D10").Value = 44215
$ws.Range("J10").Value = 16000

# Row 11
$ws.Range("D11").Value = 44230
$ws.Range("J11").Value = 16000

# Row 12
$ws.Range("D12").Value = 44168

# Row 13
$ws.Range("D13").Value = 44231
$ws.Range("J13").Value = 12000

# Row 14
$ws.Range("D14").Value = 44232
$ws.Range("J14").Value = 16000

# Row 15
$ws.Range("D15").Value = 44159

# Row 16
$ws.Range("D16").Value = 44166

# Row 18
$ws.Range("D18").Value = 44161
$ws.Range("J18").Value = 7000
$ws.Range("K18").Value = 3000
$ws.Range("M18").Value = 3000
$ws.Range("P18").Value = 30

# Row 19
$ws.Range("D19").Value = 44187
$ws.Range("J19").Value = 12000

# Row 20
$ws.Range("D20").Value = 44167
$ws.Range("J20").Value = 7000

# Row 21
$ws.Range("D21").Value = 44186
$ws.Range("J21").Value = 10000

# Row 22
$ws.Range("D22").Value = 44162
$ws.Range("K22").Value = 3000
$ws.Range("M22").Value = 3000
$ws.Range("P22").Value = 30

# Row 23
$ws.Range("D23").Value = 44189

# Row 24
$ws.Range("D24").Value = 44209
$ws.Range("K24").Value = 2500
$ws.Range("M24").Value = 2750
$ws.Range("P24").Value = 28
